# Insert a new "Status" column before column E (old E/"Saldo" shifts right to F,
# and the stray note cell that lived in F3 shifts to G3). Matches the author's
# "perbaikan data sample excel" sample-data fix: a Status column was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns E:F right by inserting a new blank column at E.
$ws.Range("E:E").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("E1").Value = "Status"

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("F9").Select()
